$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric odds values per the target diff (218 cell updates).
$updates = @{
    "Z2" = 12
    "AJ2" = 29
    "H3" = 3
    "J3" = 1.11
    "K3" = 6.5
    "N3" = 2.63
    "O3" = 1.5
    "R3" = 2.1
    "S3" = 1.67
    "V3" = 12
    "Z3" = 6.5
    "AB3" = 19
    "AI3" = 26
    "I4" = 2.9
    "J4" = 1.11
    "K4" = 6.5
    "L4" = 1.5
    "M4" = 2.63
    "N4" = 2.63
    "O4" = 1.5
    "P4" = 1.57
    "Q4" = 2.25
    "R4" = 2.05
    "S4" = 1.7
    "T4" = 7
    "X4" = 26
    "Z4" = 6.5
    "AB4" = 17
    "AC4" = 67
    "AD4" = 501
    "AE4" = 7
    "G5" = 3.3
    "P5" = 1.88
    "Q5" = 1.93
    "G7" = 2.05
    "H7" = 3.1
    "I7" = 4.2
    "J7" = 1.11
    "K7" = 6.5
    "W7" = 17
    "G8" = 3.25
    "I8" = 2.5
    "J8" = 1.14
    "K8" = 5.5
    "AA8" = 6
    "AF8" = 10
    "AH8" = 23
    "J9" = 1.11
    "K9" = 6.5
    "G13" = 4
    "H13" = 3.4
    "I13" = 1.95
    "R13" = 1.91
    "S13" = 1.8
    "U13" = 19
    "V13" = 13
    "AF13" = 8.5
    "AH13" = 17
    "G17" = 2.45
    "H17" = 3.25
    "I17" = 2.88
    "V17" = 10
    "X17" = 21
    "AA17" = 6
    "AJ17" = 29
    "G18" = 1.28
    "H18" = 5.7
    "I18" = 8.25
    "L18" = 1.09
    "M18" = 6.1
    "N18" = 1.31
    "O18" = 3.15
    "P18" = 1.19
    "Q18" = 4.15
    "R18" = 1.6
    "S18" = 2.2
    "T18" = 12.5
    "U18" = 9.25
    "W18" = 9.75
    "Y18" = 18.5
    "Z18" = 29
    "AA18" = 13
    "AB18" = 17.5
    "AC18" = 50
    "AD18" = 250
    "AE18" = 37
    "AG18" = 26
    "AH18" = 200
    "AI18" = 75
    "AJ18" = 50
    "H19" = 3.5
    "I19" = 3.05
    "R19" = 1.36
    "S19" = 2.9
    "T19" = 14.5
    "U19" = 16.5
    "W19" = 26
    "Y19" = 15.5
    "AA19" = 7.8
    "AC19" = 26
    "AE19" = 15.5
    "AF19" = 21
    "AI19" = 22
    "AJ19" = 21
    "G24" = 2.57
    "I24" = 2.55
    "T24" = 7.4
    "U24" = 12
    "AE24" = 7.6
    "AF24" = 12.5
    "AG24" = 10
    "AH24" = 28
    "H25" = 4
    "I25" = 1.57
    "J25" = 1.05
    "K25" = 11
    "Z25" = 10
    "AD25" = 451
    "AE25" = 6
    "G27" = 1.8
    "J27" = 1.07
    "K27" = 8.5
    "Z27" = 8.5
    "AD27" = 351
    "AE27" = 12
    "AF27" = 23
    "G28" = 1.36
    "H28" = 4.75
    "I28" = 8.5
    "L28" = 1.2
    "M28" = 4.33
    "N28" = 1.67
    "O28" = 2.15
    "P28" = 1.3
    "Q28" = 3.4
    "T28" = 7.5
    "U28" = 6.5
    "W28" = 9
    "X28" = 11
    "Z28" = 13
    "AA28" = 9
    "AE28" = 21
    "AG28" = 23
    "AH28" = 101
    "G34" = 1.55
    "H34" = 4
    "I34" = 4.85
    "R34" = 1.8
    "S34" = 1.91
    "T34" = 7.5
    "U34" = 7.6
    "W34" = 11.25
    "Y34" = 24
    "AA34" = 8
    "AB34" = 16.5
    "AC34" = 70
    "AD34" = 500
    "AE34" = 14.5
    "AF34" = 29
    "AG34" = 16
    "AH34" = 90
    "AI34" = 45
    "G35" = 1.42
    "H35" = 4.5
    "I35" = 7
    "N35" = 2.25
    "O35" = 1.62
    "P35" = 1.5
    "Q35" = 2.5
    "R35" = 2.63
    "S35" = 1.44
    "V35" = 9.5
    "W35" = 8.5
    "Z35" = 8
    "AA35" = 9.5
    "AB35" = 29
    "AE35" = 12
    "AF35" = 34
    "AJ35" = 67
    "G36" = 1.67
    "H36" = 3.8
    "J36" = 1.05
    "K36" = 11
    "L36" = 1.29
    "M36" = 3.5
    "N36" = 1.9
    "O36" = 1.95
    "P36" = 1.4
    "Q36" = 2.75
    "R36" = 1.91
    "S36" = 1.91
    "T36" = 7
    "Y36" = 26
    "Z36" = 11
    "AD36" = 301
    "AF36" = 23
    "G37" = 1.65
    "H37" = 3.9
    "N37" = 1.85
    "O37" = 2
    "Y37" = 23
    "Z37" = 12
    "G38" = 3.5
    "H38" = 3.5
    "I38" = 2
    "N38" = 1.67
    "O38" = 2.15
    "W38" = 41
    "X38" = 26
    "AA38" = 7
    "AF38" = 11
    "AH38" = 19
    "AJ38" = 21
    "K40" = 5.8
    "J41" = 1.05
    "K41" = 11
    "N41" = 2
    "O41" = 1.85
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
